$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A151").Value = 111474170
$ws.Range("B151").Value = 95674
$ws.Range("D151").Value = 'LC'
$ws.Range("E151").Value = 222741
$ws.Range("F151").Value = 'Finbräken'
$ws.Range("G151").Value = 'Cystopteris montana'
$ws.Range("H151").Value = '(Lam.) Desv.'
$ws.Range("A152").Value = 111475500
$ws.Range("B152").Value = 90087
$ws.Range("E152").Value = 3298
$ws.Range("F152").Value = 'Trådticka'
$ws.Range("G152").Value = 'Climacocystis borealis'
$ws.Range("H152").Value = '(Fr.) Kotl. & Pouzar'
$ws.Range("Q152").Value = 723181.0483288103
$ws.Range("R152").Value = 7544299.950535267
$ws.Range("A153").Value = 111474396
$ws.Range("B153").Value = 89423
$ws.Range("D153").Value = 'NT'
$ws.Range("E153").Value = 5432
$ws.Range("F153").Value = 'Granticka'
$ws.Range("G153").Value = 'Porodaedalea chrysoloma'
$ws.Range("H153").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q153").Value = 723316.8847442043
$ws.Range("R153").Value = 7544194.885817799
$ws.Range("A154").Value = 111473568
$ws.Range("B154").Value = 78579
$ws.Range("D154").Value = 'NT'
$ws.Range("E154").Value = 2081
$ws.Range("F154").Value = 'Skrovellav'
$ws.Range("G154").Value = 'Lobaria scrobiculata'
$ws.Range("H154").Value = '(Scop.) DC.'
$ws.Range("P154").Value = 'kurravaara ravinen, T lm'
$ws.Range("Q154").Value = 723696.4827296173
$ws.Range("R154").Value = 7544132.847775052
$ws.Range("S154").Value = 25
$ws.Range("A155").Value = 111474929
$ws.Range("B155").Value = 89423
$ws.Range("E155").Value = 5432
$ws.Range("F155").Value = 'Granticka'
$ws.Range("G155").Value = 'Porodaedalea chrysoloma'
$ws.Range("H155").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("P155").Value = 'kurra, T lm'
$ws.Range("Q155").Value = 723266.0405644479
$ws.Range("R155").Value = 7544227.158111322
$ws.Range("S155").Value = 10
$ws.Range("A156").Value = 111474891
$ws.Range("B156").Value = 90087
$ws.Range("D156").Value = 'LC'
$ws.Range("E156").Value = 3298
$ws.Range("F156").Value = 'Trådticka'
$ws.Range("G156").Value = 'Climacocystis borealis'
$ws.Range("H156").Value = '(Fr.) Kotl. & Pouzar'
$ws.Range("A157").Value = 111493799
$ws.Range("B157").Value = 78578
$ws.Range("E157").Value = 6458
$ws.Range("F157").Value = 'Lunglav'
$ws.Range("G157").Value = 'Lobaria pulmonaria'
$ws.Range("H157").Value = '(L.) Hoffm.'
$ws.Range("I157").ClearContents()
$ws.Range("J157").ClearContents()
$ws.Range("Q157").Value = 723316.8847442043
$ws.Range("R157").Value = 7544194.885817799
$ws.Range("AC157").ClearContents()
$ws.Range("AM157").Value = 'Sten/berg på land'
$ws.Range("AO157").Value = 'Stone/rock on land'
$ws.Range("A158").Value = 111494950
$ws.Range("B158").Value = 96346
$ws.Range("E158").Value = 620
$ws.Range("F158").Value = 'Skogsfru'
$ws.Range("G158").Value = 'Epipogium aphyllum'
$ws.Range("H158").Value = 'Sw.'
$ws.Range("J158").Value = 'stjälkar/strån/skott'
$ws.Range("K158").ClearContents()
$ws.Range("M158").ClearContents()
$ws.Range("P158").Value = 'kurra, T lm'
$ws.Range("Q158").Value = 723181.0483288103
$ws.Range("R158").Value = 7544299.950535267
$ws.Range("AC158").Value = '9 blommor'
$ws.Range("A159").Value = 111494318
$ws.Range("B159").Value = 56543
$ws.Range("E159").Value = 103021
$ws.Range("F159").Value = 'Talltita'
$ws.Range("G159").Value = 'Poecile montanus'
$ws.Range("H159").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("I159").Value = '2'
$ws.Range("K159").Value = 'pulli'
$ws.Range("M159").Value = 'födosökande'
$ws.Range("P159").Value = 'kurra (kurra), T lm'
$ws.Range("Q159").Value = 723122.7605886162
$ws.Range("R159").Value = 7544278.282202527
$ws.Range("AM159").ClearContents()
$ws.Range("AO159").ClearContents()
